# Issue estimates.xlsx - add new sprint planning rows + research doc progress
# to the last sprint block (rows 42-49 -> 42-50) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Formatting first: clone cell formats from pristine rows so the new /
#    shifted cells land on the SAME style indices the rest of the table
#    already uses (style 2 / 1 / 7 for columns A/B/C, style 9 for plain D
#    cells, style 10 for the SUM() formula cell).
# ---------------------------------------------------------------------------

# A:C (issue#, description, points) formatting for the new rows 46-50
$ws.Range("A45:C45").Copy()
$ws.Range("A46:C46").PasteSpecial(-4122)
$ws.Range("A47:C47").PasteSpecial(-4122)
$ws.Range("A48:C48").PasteSpecial(-4122)
$ws.Range("A49:C49").PasteSpecial(-4122)
$ws.Range("A50:C50").PasteSpecial(-4122)

# D formatting ("plain" style, as used by D43/D44/D45/D46/D49 before the edit)
$ws.Range("D43").Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("D50").PasteSpecial(-4122)

# D48 needs the SUM-formula style (same as the old D47)
$ws.Range("D47").Copy()
$ws.Range("D48").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# D42 becomes a blank "section divider" cell using the Good cell-style
# (green fill / green font) without the table border.
$ws.Range("D42").Style = "Good"
$ws.Range("D42").Borders.LineStyle = 0
$ws.Range("D42").ClearContents()

# ---------------------------------------------------------------------------
# 2) Values - write the B-column description text in the same order the
#    strings first appear so they append to the shared-string table in the
#    expected sequence.
# ---------------------------------------------------------------------------

$ws.Range("B46").Value = "Add title and description in my evidence page"
$ws.Range("B47").Value = "Add a tooltip on evidence types"
$ws.Range("B50").Value = "Add a tooltip for goal tag status"
$ws.Range("B49").Value = "Add status color on goal tags"
$ws.Range("B48").Value = "Goal and collection tags need margin"

# Days worked / Availability / Planned points labels shift down one row
$ws.Range("D43").Value = "Days worked:"
$ws.Range("D45").Value = "Planned points:"
$ws.Range("D47").Value = "Availability:"
$ws.Range("D49").Value = "Completed points:"
$ws.Range("D50").ClearContents()

# Issue numbers (column A)
$ws.Range("A46").Value = 3411
$ws.Range("A47").Value = 3410
$ws.Range("A48").Value = 3407
$ws.Range("A49").Value = 3408
$ws.Range("A50").Value = 3409

# Points spent (column C)
$ws.Range("C44").Value = 0.5
$ws.Range("C46").Value = 0.5
$ws.Range("C47").Value = 0.5
$ws.Range("C48").Value = 0.5
$ws.Range("C49").Value = 0.5
$ws.Range("C50").Value = 0.5

# Remaining numeric D-column values
$ws.Range("D44").Value = 9
$ws.Range("D46").Value = 4.5

# Total points formula, now spanning the full (bigger) block of issues
$ws.Range("D48").Formula = "=SUM(C43:C50)"

# ---------------------------------------------------------------------------
# 3) View state - mirror the author's final selection / scroll position.
# ---------------------------------------------------------------------------

$ws.Range("G45").Select()
$excel.ActiveWindow.ScrollRow = 34

Write-Host "edit applied"
